# "Generate Report for Handback" — mark the two localized files as handed
# back: zh-cn files are in sync already (no new handback timestamp shown),
# de-de files just received a fresh handback with a new timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/003644062fb98456c9c7642180c39f0846bbb010/e2e/"
$mdC5 = "c5e5ed8e-9d47-4f9a-b5f8-55cf8bb7547f.md"
$mdD4 = "d4e93497-4b60-497b-a193-6f002878760b.md"

# --- Status text, globally shared: "Ready for handoff" -> "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File / Latest Handback File populated; handback datetime
#     stays the placeholder-turned-real-timestamp string (same cell reference
#     as before, text itself changes workbook-wide).
$zhcn.Range("I2").Value = $mdC5
$zhcn.Range("J2").Value = "c5e5ed8e-9d47-4f9a-b5f8-55cf8bb7547f.4be8517b10bef7956ce91ae388fc61ff7586cb0c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-29 22:50:25"

$zhcn.Range("I3").Value = $mdD4
$zhcn.Range("J3").Value = "d4e93497-4b60-497b-a193-6f002878760b.da24cdee7f7cbd927b416865bc026b72c4954497.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-29 22:50:25"

# --- de-de: same, but with its own fresh handback datetime
$dede.Range("I2").Value = $mdC5
$dede.Range("J2").Value = "c5e5ed8e-9d47-4f9a-b5f8-55cf8bb7547f.4be8517b10bef7956ce91ae388fc61ff7586cb0c.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 22:50:33"

$dede.Range("I3").Value = $mdD4
$dede.Range("J3").Value = "d4e93497-4b60-497b-a193-6f002878760b.da24cdee7f7cbd927b416865bc026b72c4954497.de-de.xlf"
$dede.Range("K3").Value = "2016-08-29 22:50:33"

# --- New hyperlinks on the "Latest Target File" cells, pointing at the same
#     source-file links as column A, inserted so relationship ids land right
#     after each row's existing "A" hyperlink.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), ($ghBase + $mdC5), "", "", $mdC5)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($ghBase + $mdC5), "", "", $mdC5)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), ($ghBase + $mdD4), "", "", $mdD4)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($ghBase + $mdD4), "", "", $mdD4)

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), ($ghBase + $mdC5), "", "", $mdC5)
$dede.Hyperlinks.Add($dede.Range("I2"), ($ghBase + $mdC5), "", "", $mdC5)
$dede.Hyperlinks.Add($dede.Range("A3"), ($ghBase + $mdD4), "", "", $mdD4)
$dede.Hyperlinks.Add($dede.Range("I3"), ($ghBase + $mdD4), "", "", $mdD4)

# --- Column widths: widen Status + Target/Handback file columns to fit the
#     new, longer contents.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.1
$zhcn.Columns.Item(10).ColumnWidth = 39.1

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.1
$dede.Columns.Item(10).ColumnWidth = 39.1
